$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1419
$ws.Range("K3").Value = 1354
$ws.Range("J4").Value = 1792
$ws.Range("K4").Value = 295
$ws.Range("K6").Value = 1716
$ws.Range("J7").Value = 29262
$ws.Range("K7").Value = 4874

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 31
$ws.Range("K4").Value = 21
$ws.Range("J7").Value = 823
$ws.Range("K7").Value = 139
$ws.Range("J8").Value = 1852
$ws.Range("K8").Value = 287
$ws.Range("K10").Value = 27
$ws.Range("K11").Value = 100
$ws.Range("K14").Value = 31
$ws.Range("K18").Value = 39
$ws.Range("K19").Value = 127
$ws.Range("K27").Value = 58
$ws.Range("K29").Value = 225
$ws.Range("K31").Value = 55
$ws.Range("K33").Value = 197
$ws.Range("K36").Value = 54
$ws.Range("K37").Value = 165
$ws.Range("K41").Value = 51
$ws.Range("K42").Value = 167
$ws.Range("K43").Value = 47
$ws.Range("K48").Value = 54
$ws.Range("K49").Value = 31
$ws.Range("I52").Value = 595
$ws.Range("K52").Value = 132
$ws.Range("K53").Value = 72
$ws.Range("K54").Value = 87
$ws.Range("K60").Value = 37
$ws.Range("I63").Value = 196
$ws.Range("J63").Value = 89
$ws.Range("K63").Value = 14
$ws.Range("K65").Value = 127
$ws.Range("K67").Value = 191
$ws.Range("K70").Value = 10
$ws.Range("K71").Value = 14
$ws.Range("K77").Value = 36
$ws.Range("K78").Value = 71
$ws.Range("K83").Value = 98
$ws.Range("K85").Value = 249
$ws.Range("K88").Value = 62
$ws.Range("K89").Value = 66
$ws.Range("K90").Value = 46
$ws.Range("K96").Value = 67
$ws.Range("K99").Value = 90
$ws.Range("J101").Value = 29262
$ws.Range("K101").Value = 4874

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 53
$ws.Range("J4").Value = 33
$ws.Range("J7").Value = 823
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 92
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 32
$ws.Range("K3").Value = 29
$ws.Range("I5").Value = 21
$ws.Range("K6").Value = 61
$ws.Range("I7").Value = 595
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 16
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 84
$ws.Range("J4").Value = 97
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 95
$ws.Range("J7").Value = 1852
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 56
$ws.Range("K3").Value = 78
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 33
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 34
$ws.Range("K3").Value = 34
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 8
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 18
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 72
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 18
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 6
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 21
